# The document's "Requisitos" section ends with a line such as:
#   "LOB1024: Mecânica (Requisito fraco)"
# followed by an (empty) spacer paragraph, and then two site-footer
# paragraphs ("Ver no Jupiter Salvar em pdf Salvar em docx" and the
# "© 2020 . Contact: ..." copyright line) plus an extra empty spacer
# paragraph right before the trailing page-break paragraph.
#
# This rebuild of the site dropped that footer block entirely, leaving
# just the single empty paragraph that precedes the page-break
# paragraph. Locate the requirement line as an anchor and remove the
# three paragraphs that followed the existing spacer paragraph.

$d = $word.ActiveDocument

# Locate the anchor paragraph ("LOB1024: ...") by scanning paragraph text
# (robust against any shifts elsewhere in the document).
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -match "LOB1024") {
        $anchorIndex = $i
        break
    }
}

if ($anchorIndex -gt 0) {
    # anchorIndex     -> "LOB1024: Mecânica (Requisito fraco)"   (kept)
    # anchorIndex + 1 -> empty spacer paragraph                   (kept)
    # anchorIndex + 2 -> "Ver no Jupiter Salvar em pdf Salvar em docx"      (remove)
    # anchorIndex + 3 -> "© 2020 . Contact: ..." copyright line             (remove)
    # anchorIndex + 4 -> empty spacer paragraph                            (remove)
    # anchorIndex + 5 -> page-break paragraph                              (kept)
    $firstToRemove = $d.Paragraphs($anchorIndex + 2)
    $lastToRemove = $d.Paragraphs($anchorIndex + 4)

    $killRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $killRange.Delete()
}
